$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for the FID and Paq2Piq metrics (columns O and P).
# Copy the formatting (bold, centered, bordered) from the last existing
# header cell (N1) onto the two new header cells before setting their text.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("O1").Value = "FID"
$ws.Range("P1").Value = "Paq2Piq"

# New FID values per row, Paq2Piq values are all 0
$fid = @{
    2 = 29.4329585126558
    3 = 22.31222986630225
    4 = 58.43144240361824
    5 = 86.8405846731801
    6 = 92.32639214517206
    7 = 57.86872152018569
}

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 15).Value = $fid[$row]
    $ws.Cells.Item($row, 16).Value = 0
}
